$d = $word.ActiveDocument

# The document currently ends with the paragraph containing
# "This is written in new paragraph". We need to append two more
# paragraphs after it:
#   1) an empty paragraph (same paragraph formatting)
#   2) a paragraph containing a long run of "s" characters
# Word inherits paragraph/run formatting (ind firstLine=720, sz=36,
# szCs=36, lang=en-US) from the paragraph we split off of, so no
# extra formatting calls are required.

$lastPara = $d.Paragraphs.Last
$rng = $lastPara.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()

$emptyParaRng = $d.Paragraphs.Last.Range
$emptyParaRng.Collapse(0)
$emptyParaRng.InsertParagraphAfter()

$sText = "".PadRight(317, "s")
$finalParaRng = $d.Paragraphs.Last.Range
$finalParaRng.Collapse(0)
$finalParaRng.InsertAfter($sText)
